# The dataset (rows 224..262) is a weekly price log for Cilantro at
# "Terminal La Palmera de La Serena". A new week's record is inserted
# at row 224 (pushing the existing rows 224-262 down to 225-263), and
# filled in with the new week's data; all the other rows keep their
# original content (Excel's row-insert shifts them down automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 224, shifting rows 224:262 down to 225:263.
$ws.Rows.Item(224).Insert()

# Fill the new row 224 with the new weekly record. Columns A, B, C, E,
# F, G, H, I, N, O, Q, R are constant for every row in this block.
$ws.Cells.Item(224, 1).Value = 8
$ws.Cells.Item(224, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(224, 3).Value = "Coquimbo"
$ws.Cells.Item(224, 4).Value = 45218
$ws.Cells.Item(224, 5).Value = 4
$ws.Cells.Item(224, 6).Value = 100112040
$ws.Cells.Item(224, 7).Value = "Cilantro"
$ws.Cells.Item(224, 8).Value = "Sin especificar"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 2200
$ws.Cells.Item(224, 11).Value = 1500
$ws.Cells.Item(224, 12).Value = 2000
$ws.Cells.Item(224, 13).Value = 1750
$ws.Cells.Item(224, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(224, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(224, 16).Value = 1167
$ws.Cells.Item(224, 17).Value = 1.5
$ws.Cells.Item(224, 18).Value = "Hortaliza"
